# Auto update Excel log
# Appends new sensor reading rows to the PIR, Humidity, and Temperature sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("PIR")
$ws.Cells.Item(61, 1).Value = "'2026-01-28"
$ws.Cells.Item(61, 2).Value = "16:40:06"
$ws.Cells.Item(61, 3).Value = "16:00"
$ws.Cells.Item(61, 4).Value = "Bathroom"
$ws.Cells.Item(61, 5).Value = "No Motion"
$ws.Cells.Item(61, 6).Value = "Inactive"
$ws.Cells.Item(62, 1).Value = "'2026-01-28"
$ws.Cells.Item(62, 2).Value = "16:40:08"
$ws.Cells.Item(62, 3).Value = "16:00"
$ws.Cells.Item(62, 4).Value = "Bathroom"
$ws.Cells.Item(62, 5).Value = "No Motion"
$ws.Cells.Item(62, 6).Value = "Inactive"
$ws.Cells.Item(63, 1).Value = "'2026-01-28"
$ws.Cells.Item(63, 2).Value = "16:40:12"
$ws.Cells.Item(63, 3).Value = "16:00"
$ws.Cells.Item(63, 4).Value = "Bathroom"
$ws.Cells.Item(63, 5).Value = "No Motion"
$ws.Cells.Item(63, 6).Value = "Inactive"
$ws.Cells.Item(64, 1).Value = "'2026-01-28"
$ws.Cells.Item(64, 2).Value = "16:40:17"
$ws.Cells.Item(64, 3).Value = "16:00"
$ws.Cells.Item(64, 4).Value = "Bathroom"
$ws.Cells.Item(64, 5).Value = "No Motion"
$ws.Cells.Item(64, 6).Value = "Inactive"
$ws.Cells.Item(65, 1).Value = "'2026-01-28"
$ws.Cells.Item(65, 2).Value = "16:40:23"
$ws.Cells.Item(65, 3).Value = "16:00"
$ws.Cells.Item(65, 4).Value = "Bathroom"
$ws.Cells.Item(65, 5).Value = "No Motion"
$ws.Cells.Item(65, 6).Value = "Inactive"
$ws.Cells.Item(66, 1).Value = "'2026-01-28"
$ws.Cells.Item(66, 2).Value = "16:40:28"
$ws.Cells.Item(66, 3).Value = "16:00"
$ws.Cells.Item(66, 4).Value = "Bathroom"
$ws.Cells.Item(66, 5).Value = "No Motion"
$ws.Cells.Item(66, 6).Value = "Inactive"
$ws.Cells.Item(67, 1).Value = "'2026-01-28"
$ws.Cells.Item(67, 2).Value = "16:40:33"
$ws.Cells.Item(67, 3).Value = "16:00"
$ws.Cells.Item(67, 4).Value = "Bathroom"
$ws.Cells.Item(67, 5).Value = "No Motion"
$ws.Cells.Item(67, 6).Value = "Inactive"
$ws.Cells.Item(68, 1).Value = "'2026-01-28"
$ws.Cells.Item(68, 2).Value = "16:40:38"
$ws.Cells.Item(68, 3).Value = "16:00"
$ws.Cells.Item(68, 4).Value = "Bathroom"
$ws.Cells.Item(68, 5).Value = "No Motion"
$ws.Cells.Item(68, 6).Value = "Inactive"
$ws.Cells.Item(69, 1).Value = "'2026-01-28"
$ws.Cells.Item(69, 2).Value = "16:40:43"
$ws.Cells.Item(69, 3).Value = "16:00"
$ws.Cells.Item(69, 4).Value = "Bathroom"
$ws.Cells.Item(69, 5).Value = "No Motion"
$ws.Cells.Item(69, 6).Value = "Inactive"
$ws.Cells.Item(70, 1).Value = "'2026-01-28"
$ws.Cells.Item(70, 2).Value = "16:40:48"
$ws.Cells.Item(70, 3).Value = "16:00"
$ws.Cells.Item(70, 4).Value = "Bathroom"
$ws.Cells.Item(70, 5).Value = "No Motion"
$ws.Cells.Item(70, 6).Value = "Inactive"
$ws.Cells.Item(71, 1).Value = "'2026-01-28"
$ws.Cells.Item(71, 2).Value = "16:40:53"
$ws.Cells.Item(71, 3).Value = "16:00"
$ws.Cells.Item(71, 4).Value = "Bathroom"
$ws.Cells.Item(71, 5).Value = "No Motion"
$ws.Cells.Item(71, 6).Value = "Inactive"
$ws.Cells.Item(72, 1).Value = "'2026-01-28"
$ws.Cells.Item(72, 2).Value = "16:40:58"
$ws.Cells.Item(72, 3).Value = "16:00"
$ws.Cells.Item(72, 4).Value = "Bathroom"
$ws.Cells.Item(72, 5).Value = "No Motion"
$ws.Cells.Item(72, 6).Value = "Inactive"
$ws.Cells.Item(73, 1).Value = "'2026-01-28"
$ws.Cells.Item(73, 2).Value = "16:41:03"
$ws.Cells.Item(73, 3).Value = "16:00"
$ws.Cells.Item(73, 4).Value = "Bathroom"
$ws.Cells.Item(73, 5).Value = "No Motion"
$ws.Cells.Item(73, 6).Value = "Inactive"

$ws = $wb.Worksheets.Item("Humidity")
$ws.Cells.Item(60, 1).Value = "'2026-01-28"
$ws.Cells.Item(60, 2).Value = "16:40:07"
$ws.Cells.Item(60, 3).Value = "16:00"
$ws.Cells.Item(60, 4).Value = "Bathroom"
$ws.Cells.Item(60, 5).Value = "'87.9%"
$ws.Cells.Item(60, 6).Value = "Active"
$ws.Cells.Item(61, 1).Value = "'2026-01-28"
$ws.Cells.Item(61, 2).Value = "16:40:11"
$ws.Cells.Item(61, 3).Value = "16:00"
$ws.Cells.Item(61, 4).Value = "Bathroom"
$ws.Cells.Item(61, 5).Value = "'87.9%"
$ws.Cells.Item(61, 6).Value = "Active"
$ws.Cells.Item(62, 1).Value = "'2026-01-28"
$ws.Cells.Item(62, 2).Value = "16:40:15"
$ws.Cells.Item(62, 3).Value = "16:00"
$ws.Cells.Item(62, 4).Value = "Bathroom"
$ws.Cells.Item(62, 5).Value = "'87.0%"
$ws.Cells.Item(62, 6).Value = "Active"
$ws.Cells.Item(63, 1).Value = "'2026-01-28"
$ws.Cells.Item(63, 2).Value = "16:40:19"
$ws.Cells.Item(63, 3).Value = "16:00"
$ws.Cells.Item(63, 4).Value = "Bathroom"
$ws.Cells.Item(63, 5).Value = "'87.9%"
$ws.Cells.Item(63, 6).Value = "Active"
$ws.Cells.Item(64, 1).Value = "'2026-01-28"
$ws.Cells.Item(64, 2).Value = "16:40:23"
$ws.Cells.Item(64, 3).Value = "16:00"
$ws.Cells.Item(64, 4).Value = "Bathroom"
$ws.Cells.Item(64, 5).Value = "'87.9%"
$ws.Cells.Item(64, 6).Value = "Active"
$ws.Cells.Item(65, 1).Value = "'2026-01-28"
$ws.Cells.Item(65, 2).Value = "16:40:27"
$ws.Cells.Item(65, 3).Value = "16:00"
$ws.Cells.Item(65, 4).Value = "Bathroom"
$ws.Cells.Item(65, 5).Value = "'87.0%"
$ws.Cells.Item(65, 6).Value = "Active"
$ws.Cells.Item(66, 1).Value = "'2026-01-28"
$ws.Cells.Item(66, 2).Value = "16:40:31"
$ws.Cells.Item(66, 3).Value = "16:00"
$ws.Cells.Item(66, 4).Value = "Bathroom"
$ws.Cells.Item(66, 5).Value = "'87.9%"
$ws.Cells.Item(66, 6).Value = "Active"
$ws.Cells.Item(67, 1).Value = "'2026-01-28"
$ws.Cells.Item(67, 2).Value = "16:40:35"
$ws.Cells.Item(67, 3).Value = "16:00"
$ws.Cells.Item(67, 4).Value = "Bathroom"
$ws.Cells.Item(67, 5).Value = "'87.0%"
$ws.Cells.Item(67, 6).Value = "Active"
$ws.Cells.Item(68, 1).Value = "'2026-01-28"
$ws.Cells.Item(68, 2).Value = "16:40:43"
$ws.Cells.Item(68, 3).Value = "16:00"
$ws.Cells.Item(68, 4).Value = "Bathroom"
$ws.Cells.Item(68, 5).Value = "'87.9%"
$ws.Cells.Item(68, 6).Value = "Active"
$ws.Cells.Item(69, 1).Value = "'2026-01-28"
$ws.Cells.Item(69, 2).Value = "16:40:47"
$ws.Cells.Item(69, 3).Value = "16:00"
$ws.Cells.Item(69, 4).Value = "Bathroom"
$ws.Cells.Item(69, 5).Value = "'87.0%"
$ws.Cells.Item(69, 6).Value = "Active"
$ws.Cells.Item(70, 1).Value = "'2026-01-28"
$ws.Cells.Item(70, 2).Value = "16:40:51"
$ws.Cells.Item(70, 3).Value = "16:00"
$ws.Cells.Item(70, 4).Value = "Bathroom"
$ws.Cells.Item(70, 5).Value = "'87.9%"
$ws.Cells.Item(70, 6).Value = "Active"
$ws.Cells.Item(71, 1).Value = "'2026-01-28"
$ws.Cells.Item(71, 2).Value = "16:40:59"
$ws.Cells.Item(71, 3).Value = "16:00"
$ws.Cells.Item(71, 4).Value = "Bathroom"
$ws.Cells.Item(71, 5).Value = "'87.0%"
$ws.Cells.Item(71, 6).Value = "Active"
$ws.Cells.Item(72, 1).Value = "'2026-01-28"
$ws.Cells.Item(72, 2).Value = "16:41:03"
$ws.Cells.Item(72, 3).Value = "16:00"
$ws.Cells.Item(72, 4).Value = "Bathroom"
$ws.Cells.Item(72, 5).Value = "'88.0%"
$ws.Cells.Item(72, 6).Value = "Active"

$ws = $wb.Worksheets.Item("Temperature")
$ws.Cells.Item(60, 1).Value = "'2026-01-28"
$ws.Cells.Item(60, 2).Value = "16:40:07"
$ws.Cells.Item(60, 3).Value = "16:00"
$ws.Cells.Item(60, 4).Value = "Bathroom"
$ws.Cells.Item(60, 5).Value = "22.8C"
$ws.Cells.Item(60, 6).Value = "Active"
$ws.Cells.Item(61, 1).Value = "'2026-01-28"
$ws.Cells.Item(61, 2).Value = "16:40:11"
$ws.Cells.Item(61, 3).Value = "16:00"
$ws.Cells.Item(61, 4).Value = "Bathroom"
$ws.Cells.Item(61, 5).Value = "22.8C"
$ws.Cells.Item(61, 6).Value = "Active"
$ws.Cells.Item(62, 1).Value = "'2026-01-28"
$ws.Cells.Item(62, 2).Value = "16:40:15"
$ws.Cells.Item(62, 3).Value = "16:00"
$ws.Cells.Item(62, 4).Value = "Bathroom"
$ws.Cells.Item(62, 5).Value = "22.8C"
$ws.Cells.Item(62, 6).Value = "Active"
$ws.Cells.Item(63, 1).Value = "'2026-01-28"
$ws.Cells.Item(63, 2).Value = "16:40:19"
$ws.Cells.Item(63, 3).Value = "16:00"
$ws.Cells.Item(63, 4).Value = "Bathroom"
$ws.Cells.Item(63, 5).Value = "22.8C"
$ws.Cells.Item(63, 6).Value = "Active"
$ws.Cells.Item(64, 1).Value = "'2026-01-28"
$ws.Cells.Item(64, 2).Value = "16:40:23"
$ws.Cells.Item(64, 3).Value = "16:00"
$ws.Cells.Item(64, 4).Value = "Bathroom"
$ws.Cells.Item(64, 5).Value = "22.8C"
$ws.Cells.Item(64, 6).Value = "Active"
$ws.Cells.Item(65, 1).Value = "'2026-01-28"
$ws.Cells.Item(65, 2).Value = "16:40:27"
$ws.Cells.Item(65, 3).Value = "16:00"
$ws.Cells.Item(65, 4).Value = "Bathroom"
$ws.Cells.Item(65, 5).Value = "22.8C"
$ws.Cells.Item(65, 6).Value = "Active"
$ws.Cells.Item(66, 1).Value = "'2026-01-28"
$ws.Cells.Item(66, 2).Value = "16:40:31"
$ws.Cells.Item(66, 3).Value = "16:00"
$ws.Cells.Item(66, 4).Value = "Bathroom"
$ws.Cells.Item(66, 5).Value = "22.8C"
$ws.Cells.Item(66, 6).Value = "Active"
$ws.Cells.Item(67, 1).Value = "'2026-01-28"
$ws.Cells.Item(67, 2).Value = "16:40:36"
$ws.Cells.Item(67, 3).Value = "16:00"
$ws.Cells.Item(67, 4).Value = "Bathroom"
$ws.Cells.Item(67, 5).Value = "22.8C"
$ws.Cells.Item(67, 6).Value = "Active"
$ws.Cells.Item(68, 1).Value = "'2026-01-28"
$ws.Cells.Item(68, 2).Value = "16:40:43"
$ws.Cells.Item(68, 3).Value = "16:00"
$ws.Cells.Item(68, 4).Value = "Bathroom"
$ws.Cells.Item(68, 5).Value = "22.8C"
$ws.Cells.Item(68, 6).Value = "Active"
$ws.Cells.Item(69, 1).Value = "'2026-01-28"
$ws.Cells.Item(69, 2).Value = "16:40:48"
$ws.Cells.Item(69, 3).Value = "16:00"
$ws.Cells.Item(69, 4).Value = "Bathroom"
$ws.Cells.Item(69, 5).Value = "22.8C"
$ws.Cells.Item(69, 6).Value = "Active"
$ws.Cells.Item(70, 1).Value = "'2026-01-28"
$ws.Cells.Item(70, 2).Value = "16:40:52"
$ws.Cells.Item(70, 3).Value = "16:00"
$ws.Cells.Item(70, 4).Value = "Bathroom"
$ws.Cells.Item(70, 5).Value = "22.8C"
$ws.Cells.Item(70, 6).Value = "Active"
$ws.Cells.Item(71, 1).Value = "'2026-01-28"
$ws.Cells.Item(71, 2).Value = "16:41:00"
$ws.Cells.Item(71, 3).Value = "16:00"
$ws.Cells.Item(71, 4).Value = "Bathroom"
$ws.Cells.Item(71, 5).Value = "22.8C"
$ws.Cells.Item(71, 6).Value = "Active"
$ws.Cells.Item(72, 1).Value = "'2026-01-28"
$ws.Cells.Item(72, 2).Value = "16:41:04"
$ws.Cells.Item(72, 3).Value = "16:00"
$ws.Cells.Item(72, 4).Value = "Bathroom"
$ws.Cells.Item(72, 5).Value = "22.8C"
$ws.Cells.Item(72, 6).Value = "Active"
